$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("2025")
$ws1.Range("B2").Value = 973.9537847600009
$ws1.Range("E2").Value = 28982.37596598056
$ws1.Range("I2").Value = 16175.28135478
$ws1.Range("L2").Value = 48524.529503538
$ws1.Range("M2").Value = 10590.587968015
$ws1.Range("N2").Value = 7098.474579830177
$ws1.Range("O2").Value = 6965.142708907404

$ws2 = $wb.Worksheets.Item("2030")
$ws2.Range("B2").Value = 5712.560177842886
$ws2.Range("E2").Value = 56106.05588781912
$ws2.Range("I2").Value = 44217.8984721661
$ws2.Range("L2").Value = 66966.57749858923
$ws2.Range("M2").Value = 21984.28023276101
$ws2.Range("N2").Value = 10503.58911838769
$ws2.Range("O2").Value = 12022.61948454972

$ws3 = $wb.Worksheets.Item("2035")
$ws3.Range("A2").Value = 2861.961401238371
$ws3.Range("B2").Value = 8026.889663087295
$ws3.Range("E2").Value = 67297.73995507321
$ws3.Range("I2").Value = 59256.42575923612
$ws3.Range("L2").Value = 66966.57749858923
$ws3.Range("M2").Value = 25464.6214365565
$ws3.Range("N2").Value = 15008.88598672954
$ws3.Range("O2").Value = 14717.53372517132

$ws4 = $wb.Worksheets.Item("2040")
$ws4.Range("A2").Value = 2861.961401238371
$ws4.Range("B2").Value = 8026.889663087295
$ws4.Range("E2").Value = 67297.73995507321
$ws4.Range("I2").Value = 59256.42575923612
$ws4.Range("L2").Value = 66966.57749858923
$ws4.Range("M2").Value = 25464.6214365565
$ws4.Range("N2").Value = 15111.10479697713
$ws4.Range("O2").Value = 14717.53372517132

$ws5 = $wb.Worksheets.Item("2045")
$ws5.Range("A2").Value = 6302.873118834019
$ws5.Range("B2").Value = 8026.889663087295
$ws5.Range("E2").Value = 67297.73995507321
$ws5.Range("I2").Value = 59256.42575923612
$ws5.Range("L2").Value = 66966.57749858923
$ws5.Range("M2").Value = 25464.6214365565
$ws5.Range("N2").Value = 15633.82115826023
$ws5.Range("O2").Value = 17024.64289623143

$ws6 = $wb.Worksheets.Item("2050")
$ws6.Range("A2").Value = 6302.873118834019
$ws6.Range("B2").Value = 8026.889663087295
$ws6.Range("E2").Value = 67297.73995507321
$ws6.Range("I2").Value = 59256.42575923612
$ws6.Range("L2").Value = 66966.57749858923
$ws6.Range("M2").Value = 25464.6214365565
$ws6.Range("N2").Value = 15633.82115826023
$ws6.Range("O2").Value = 17024.64289623143
